$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply updated odds values scraped for 2025-02-01 matches
$updates = @(
    @(10, 7, 1.5),
    @(10, 9, 6.5),
    @(10, 10, 2.05),
    @(10, 12, 6),
    @(10, 15, 1.22),
    @(10, 16, 4.33),
    @(10, 17, 1.73),
    @(10, 18, 2.1),
    @(10, 19, 2.75),
    @(10, 20, 1.44),
    @(10, 25, 7.5),
    @(10, 26, 7.5),
    @(10, 35, 251),
    @(10, 40, 41),
    @(10, 41, 41),
    @(20, 9, 3.9),
    @(20, 12, 4.33),
    @(20, 17, 2),
    @(20, 18, 1.8),
    @(20, 23, 1.8),
    @(20, 24, 1.95),
    @(20, 31, 9.5),
    @(20, 35, 251),
    @(23, 8, 4.2),
    @(23, 9, 1.45),
    @(23, 12, 1.95),
    @(23, 17, 1.48),
    @(23, 18, 2.6),
    @(23, 19, 2.1),
    @(23, 20, 1.67),
    @(23, 27, 19),
    @(23, 31, 19),
    @(23, 35, 126),
    @(23, 37, 9),
    @(23, 41, 19),
    @(27, 10, 7.5),
    @(27, 11, 2.75),
    @(27, 15, 1.14),
    @(27, 16, 5.5),
    @(27, 17, 1.5),
    @(27, 18, 2.63),
    @(27, 19, 2.2),
    @(27, 20, 1.67),
    @(27, 21, 1.25),
    @(27, 22, 3.75),
    @(27, 23, 1.8),
    @(27, 24, 1.95),
    @(27, 25, 26),
    @(27, 26, 51),
    @(27, 27, 26),
    @(27, 32, 10),
    @(27, 33, 19),
    @(27, 35, 201),
    @(27, 36, 9),
    @(27, 37, 7.5),
    @(27, 41, 23),
    @(27, 42, 1.83),
    @(27, 43, 2.07),
    @(28, 8, 2.75),
    @(28, 9, 2.8),
    @(28, 11, 1.8),
    @(28, 12, 3.75),
    @(28, 13, 1.17),
    @(28, 14, 5),
    @(28, 15, 1.73),
    @(28, 16, 2),
    @(28, 17, 3.4),
    @(28, 18, 1.33),
    @(28, 19, 7),
    @(28, 20, 1.1),
    @(28, 21, 1.75),
    @(28, 22, 2.05),
    @(28, 23, 2.5),
    @(28, 24, 1.5),
    @(28, 25, 6),
    @(28, 26, 12),
    @(28, 31, 5),
    @(28, 34, 101),
    @(28, 40, 34),
    @(28, 41, 51),
    @(28, 42, 6.2),
    @(28, 43, 1.13),
    @(29, 7, 3.4),
    @(29, 8, 2.8),
    @(29, 11, 1.8),
    @(29, 13, 1.14),
    @(29, 14, 5.5),
    @(29, 17, 3.4),
    @(29, 18, 1.33),
    @(29, 21, 1.73),
    @(29, 22, 2.08),
    @(29, 27, 15),
    @(29, 33, 23),
    @(29, 37, 9.5),
    @(29, 41, 51),
    @(29, 43, 1.15),
    @(30, 13, 1.07),
    @(30, 15, 1.47),
    @(30, 18, 1.47),
    @(30, 20, 1.13),
    @(30, 44, 1.87),
    @(30, 45, 1.87),
    @(39, 7, 1.18),
    @(39, 8, 6.5),
    @(39, 9, 13),
    @(39, 10, 1.57),
    @(39, 11, 3),
    @(39, 12, 10),
    @(39, 14, 21),
    @(39, 15, 1.11),
    @(39, 16, 6.5),
    @(39, 17, 1.4),
    @(39, 18, 2.88),
    @(39, 23, 2),
    @(39, 24, 1.75),
    @(39, 29, 11),
    @(39, 32, 13),
    @(39, 34, 67),
    @(39, 36, 34),
    @(39, 37, 67),
    @(39, 40, 81),
    @(39, 41, 67),
    @(43, 7, 2),
    @(43, 9, 4),
    @(43, 10, 2.75),
    @(43, 12, 4.75),
    @(43, 21, 1.53),
    @(43, 22, 2.38),
    @(43, 26, 8.5),
    @(43, 28, 17),
    @(43, 29, 19),
    @(43, 31, 7),
    @(43, 36, 9),
    @(43, 37, 19),
    @(43, 40, 41),
    @(43, 41, 51),
    @(44, 7, 3.75),
    @(44, 9, 2.15),
    @(44, 10, 4.5),
    @(44, 11, 1.83),
    @(44, 12, 3),
    @(44, 13, 1.13),
    @(44, 14, 6),
    @(44, 15, 1.57),
    @(44, 16, 2.25),
    @(44, 17, 2.7),
    @(44, 18, 1.44),
    @(44, 19, 6),
    @(44, 20, 1.13),
    @(44, 21, 1.62),
    @(44, 22, 2.2),
    @(44, 23, 2.25),
    @(44, 24, 1.57),
    @(44, 26, 17),
    @(44, 27, 15),
    @(44, 28, 41),
    @(44, 29, 41),
    @(44, 30, 51),
    @(44, 31, 6),
    @(44, 33, 21),
    @(44, 34, 81),
    @(44, 36, 5.5),
    @(44, 37, 8.5),
    @(44, 42, 4.4),
    @(44, 43, 1.2),
    @(44, 44, 2.1),
    @(44, 45, 1.78),
    @(45, 7, 1.7),
    @(45, 9, 5),
    @(45, 18, 1.63),
    @(45, 28, 13),
    @(45, 36, 11),
    @(45, 37, 23),
    @(112, 13, 1.04),
    @(112, 14, 13),
    @(112, 15, 1.2),
    @(112, 16, 4.33),
    @(112, 17, 1.67),
    @(112, 18, 2.15),
    @(112, 27, 8.5),
    @(112, 33, 19),
    @(112, 36, 23),
    @(112, 38, 26),
    @(113, 7, 1.45),
    @(113, 9, 7),
    @(113, 10, 2),
    @(113, 12, 6.5),
    @(113, 13, 1.04),
    @(113, 14, 13),
    @(113, 17, 1.8),
    @(113, 18, 2),
    @(113, 23, 1.95),
    @(113, 24, 1.8),
    @(113, 25, 6.5),
    @(113, 26, 7),
    @(113, 28, 10),
    @(113, 30, 26),
    @(113, 32, 8),
    @(113, 34, 51),
    @(113, 37, 34),
    @(115, 15, 1.22),
    @(115, 16, 4),
    @(115, 17, 1.7),
    @(115, 18, 2.1),
    @(115, 19, 2.75),
    @(115, 20, 1.4),
    @(126, 7, 1.55),
    @(126, 8, 3.8),
    @(126, 9, 6.25),
    @(126, 10, 2.2),
    @(126, 12, 6),
    @(126, 23, 2),
    @(126, 24, 1.75),
    @(126, 25, 6.5),
    @(126, 26, 7),
    @(126, 28, 11),
    @(126, 32, 7),
    @(126, 33, 19),
    @(126, 35, 401),
    @(126, 37, 29),
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
